$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Template row used to copy cell formatting (styles) for header separator rows
$headerTemplateRange = $ws.Range("A1097:F1097")

# Row 1110: header separator "WEDNESDAY"
$headerTemplateRange.Copy()
$ws.Range("A1110:F1110").PasteSpecial(-4122)
$ws.Range("B1110").Value = 'WEDNESDAY'

# Row 1111
$ws.Range("A1111").Value = 'Pickup Mic'
$ws.Range("B1111").Value = 42753
$ws.Range("C1111").Value = '1700'
$ws.Range("D1111").Value = 'SC'
$ws.Range("E1111").Value = 'MDR'
$ws.Range("F1111").Value = 'Pick up Lecturn mic and stand and cable.  Return to Bethune 201 storeroom.'
$ws.Rows.Item(1111).RowHeight = 30

# Row 1112
$ws.Range("A1112").Value = 'Pickup Small PA'
$ws.Range("B1112").Value = 42753
$ws.Range("C1112").Value = '1700'
$ws.Range("D1112").Value = 'SC'
$ws.Range("E1112").Value = 'MDR'
$ws.Range("F1112").Value = 'Pick up speaker and cart and matts and ac cords. Return to Bethune 201 storeroom. Key for room in CB 121A storeroom.'
$ws.Rows.Item(1112).RowHeight = 30

# Row 1113
$ws.Range("A1113").Value = 'Pickup PC'
$ws.Range("B1113").Value = 42753
$ws.Range("C1113").Value = '1700'
$ws.Range("D1113").Value = 'SC'
$ws.Range("E1113").Value = 'MDR'
$ws.Range("F1113").Value = 'Pick up roll in PC and Projector carts, all matts and cables and return equipment to Bethune 201 storeroom. Key for Stong MDR is in CB 121A storeroom.'
$ws.Rows.Item(1113).RowHeight = 45

# Row 1114
$ws.Range("A1114").Value = 'Other'
$ws.Range("B1114").Value = 42753
$ws.Range("C1114").Value = '1700'
$ws.Range("D1114").Value = 'SC'
$ws.Range("E1114").Value = 'MDR'
$ws.Range("F1114").Value = 'LEAVE PORTABLE SCREEN IN ROOM'

# Row 1115
$ws.Range("A1115").Value = 'Setup Mic'
$ws.Range("B1115").Value = 42753
$ws.Range("C1115").Value = '1715'
$ws.Range("D1115").Value = 'LAS'
$ws.Range("E1115").Value = 'C'
$ws.Range("F1115").Value = 'Take cart with mixer, 2 wireless mics and 2 mic stands from Lassonde 1011 storeroom (across from Lassonde A). Go to Lassonde C classroom (class starts at 5:30 pm but be there early in case previous class ends early). '
$ws.Rows.Item(1115).RowHeight = 60

# Row 1116
$ws.Range("A1116").Value = 'Setup Mic'
$ws.Range("B1116").Value = 42753
$ws.Range("C1116").Value = '1715'
$ws.Range("D1116").Value = 'LAS'
$ws.Range("E1116").Value = 'C'
$ws.Range("F1116").Value = 'Log in as 5065*0 on touchscreen. (First level bar is your wireless handheld mic volume). Plug in mic cable from output of mixer to mic input on podium (XLR jack just above VHS machine in podium). Ramp up volume a bit on "Microphone 2" on touchscreen to medium volume to get level.'
$ws.Rows.Item(1116).RowHeight = 75

# Row 1117
$ws.Range("A1117").Value = 'Setup Mic'
$ws.Range("B1117").Value = 42753
$ws.Range("C1117").Value = '1715'
$ws.Range("D1117").Value = 'LAS'
$ws.Range("E1117").Value = 'C'
$ws.Range("F1117").Value = 'Plug in power cord from cart on to power outlet on left side of podium (to left of document camera). Turn on mixer. Turn on wireless microphone receivers on cart (NOTE: DO NOT PRESS "SYNC" BUTTON" - POWER BUTTON IS FIRST BUTTON TO THE RIGHT ON RECEIVER). '
$ws.Rows.Item(1117).RowHeight = 75

# Row 1118
$ws.Range("A1118").Value = 'Setup Mic'
$ws.Range("B1118").Value = 42753
$ws.Range("C1118").Value = '1715'
$ws.Range("D1118").Value = 'LAS'
$ws.Range("E1118").Value = 'C'
$ws.Range("F1118").Value = 'Press "MUTE" button on wireless mics to turn on mics.'

# Row 1119
$ws.Range("A1119").Value = 'Setup Mic'
$ws.Range("B1119").Value = 42753
$ws.Range("C1119").Value = '1715'
$ws.Range("D1119").Value = 'LAS'
$ws.Range("E1119").Value = 'C'
$ws.Range("F1119").Value = 'Once volumes are set, place one mic stand with mic halfway up aisle on right and one mic stand with mic halfway up aisle on left. Demo volume controls to prof. and demo PC. Leave microphone bags with milk carton on cart in room. PLEASE FIND OUT END TIME OF CLASS FROM PROF. AND TELL MASI AS MICROPHONES ARE EXPENSIVE. TELL PROF. TO STAY WITH MICS UNTIL THEY ARE PICKED UP. TELL HIM TO CALL ext 55800   WHEN DONE (use phone in classroom).'
$ws.Rows.Item(1119).RowHeight = 120

# Row 1120
$ws.Range("A1120").Value = 'AV Shutdown'
$ws.Range("B1120").Value = 42753
$ws.Range("C1120").Value = '1830'
$ws.Range("D1120").Value = 'CLH'
$ws.Range("E1120").Value = 'I'
$ws.Range("F1120").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1121
$ws.Range("A1121").Value = 'AV Shutdown'
$ws.Range("B1121").Value = 42753
$ws.Range("C1121").Value = '1900'
$ws.Range("D1121").Value = 'CLH'
$ws.Range("E1121").Value = 'L'
$ws.Range("F1121").Value = 'PLEASE MAKE SURE CRESTRON GETS LOGGED OFF. WE ARE HAVING PROBLEMS WITH THIS ROOM WHEN IT DOESN''T GET LOGGED OFF.'
$ws.Rows.Item(1121).RowHeight = 45

# Row 1122
$ws.Range("A1122").Value = 'AV Shutdown'
$ws.Range("B1122").Value = 42753
$ws.Range("C1122").Value = '1730'
$ws.Range("D1122").Value = 'LSB'
$ws.Range("E1122").Value = '101'
$ws.Range("F1122").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1123
$ws.Range("A1123").Value = 'AV Shutdown'
$ws.Range("B1123").Value = 42753
$ws.Range("C1123").Value = '1900'
$ws.Range("D1123").Value = 'LSB'
$ws.Range("E1123").Value = '103'
$ws.Range("F1123").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1124
$ws.Range("A1124").Value = 'AV Shutdown'
$ws.Range("B1124").Value = 42753
$ws.Range("C1124").Value = '1900'
$ws.Range("D1124").Value = 'LSB'
$ws.Range("E1124").Value = '105'
$ws.Range("F1124").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1125
$ws.Range("A1125").Value = 'AV Shutdown'
$ws.Range("B1125").Value = 42753
$ws.Range("C1125").Value = '2000'
$ws.Range("D1125").Value = 'LSB'
$ws.Range("E1125").Value = '106'
$ws.Range("F1125").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1126
$ws.Range("A1126").Value = 'AV Shutdown'
$ws.Range("B1126").Value = 42753
$ws.Range("C1126").Value = '1730'
$ws.Range("D1126").Value = 'LSB'
$ws.Range("E1126").Value = '107'
$ws.Range("F1126").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1127
$ws.Range("A1127").Value = 'Pickup Mic'
$ws.Range("B1127").Value = 42753
$ws.Range("C1127").Value = '1850'
$ws.Range("D1127").Value = 'LAS'
$ws.Range("E1127").Value = 'C'
$ws.Range("F1127").Value = 'Pick up 2 wireless mics on stands with cart. Move all equipment on cart - cart has 2 wireless mic receivers and mixer and mic cables. Pick up 2 mic stands - return all equipment to Lassonde 1011 storeroom (across the hall from Lassonde A). PLEASE PUT 2 WIRELESS MICS IN BAGS PROVIDED IN MILK CARTON ON CART. Very expensive mics - please go early and treat mics with care.'
$ws.Rows.Item(1127).RowHeight = 90

# Row 1128
$ws.Range("A1128").Value = 'Pickup Mic'
$ws.Range("B1128").Value = 42753
$ws.Range("C1128").Value = '1850'
$ws.Range("D1128").Value = 'LAS'
$ws.Range("E1128").Value = 'C'
$ws.Range("F1128").Value = 'Turn off wireless microphones by pressing "MUTE" button on mics.'
$ws.Rows.Item(1128).RowHeight = 30

# Row 1129
$ws.Range("A1129").Value = 'Pickup Mic'
$ws.Range("B1129").Value = 42753
$ws.Range("C1129").Value = '1850'
$ws.Range("D1129").Value = 'LAS'
$ws.Range("E1129").Value = 'C'
$ws.Range("F1129").Value = 'Turn off wireless microphone receivers by pressing "POWER" button and not "SYNC" button. '
$ws.Rows.Item(1129).RowHeight = 30

# Row 1130
$ws.Range("A1130").Value = 'Pickup Mic'
$ws.Range("B1130").Value = 42753
$ws.Range("C1130").Value = '1850'
$ws.Range("D1130").Value = 'LAS'
$ws.Range("E1130").Value = 'C'
$ws.Range("F1130").Value = 'PLEASE BE ON TIME - Prof upset last week when no one came till 7:05 pm and other class was starting.'
$ws.Rows.Item(1130).RowHeight = 30

# Row 1134: header separator "THURSDAY"
$headerTemplateRange.Copy()
$ws.Range("A1134:F1134").PasteSpecial(-4122)
$ws.Range("B1134").Value = 'THURSDAY'

# Row 1135
$ws.Range("A1135").Value = 'AV Shutdown'
$ws.Range("B1135").Value = 42754
$ws.Range("C1135").Value = '1730'
$ws.Range("D1135").Value = 'CLH'
$ws.Range("E1135").Value = 'L'
$ws.Range("F1135").Value = 'PLEASE MAKE SURE CRESTRON GETS LOGGED OFF. WE ARE HAVING PROBLEMS WITH THIS ROOM WHEN IT DOESN''T GET LOGGED OFF.'
$ws.Rows.Item(1135).RowHeight = 45

# Row 1136
$ws.Range("A1136").Value = 'AV Shutdown'
$ws.Range("B1136").Value = 42754
$ws.Range("C1136").Value = '1730'
$ws.Range("D1136").Value = 'LSB'
$ws.Range("E1136").Value = '101'
$ws.Range("F1136").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1137
$ws.Range("A1137").Value = 'AV Shutdown'
$ws.Range("B1137").Value = 42754
$ws.Range("C1137").Value = '1730'
$ws.Range("D1137").Value = 'LSB'
$ws.Range("E1137").Value = '103'
$ws.Range("F1137").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1138
$ws.Range("A1138").Value = 'AV Shutdown'
$ws.Range("B1138").Value = 42754
$ws.Range("C1138").Value = '1730'
$ws.Range("D1138").Value = 'LSB'
$ws.Range("E1138").Value = '106'
$ws.Range("F1138").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1139
$ws.Range("A1139").Value = 'AV Shutdown'
$ws.Range("B1139").Value = 42754
$ws.Range("C1139").Value = '1730'
$ws.Range("D1139").Value = 'LSB'
$ws.Range("E1139").Value = '107'
$ws.Range("F1139").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1140
$ws.Range("A1140").Value = 'AV Shutdown'
$ws.Range("B1140").Value = 42754
$ws.Range("C1140").Value = '1900'
$ws.Range("D1140").Value = 'LAS'
$ws.Range("E1140").Value = 'A'
$ws.Range("F1140").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1141
$ws.Range("A1141").Value = 'AV Shutdown'
$ws.Range("B1141").Value = 42754
$ws.Range("C1141").Value = '1900'
$ws.Range("D1141").Value = 'LAS'
$ws.Range("E1141").Value = 'B'
$ws.Range("F1141").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1142
$ws.Range("A1142").Value = 'AV Shutdown'
$ws.Range("B1142").Value = 42754
$ws.Range("C1142").Value = '2100'
$ws.Range("D1142").Value = 'CLH'
$ws.Range("E1142").Value = 'I'
$ws.Range("F1142").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1146: header separator "FRIDAY"
$headerTemplateRange.Copy()
$ws.Range("A1146:F1146").PasteSpecial(-4122)
$ws.Range("B1146").Value = 'FRIDAY'

# Row 1147
$ws.Range("A1147").Value = 'AV Shutdown'
$ws.Range("B1147").Value = 42755
$ws.Range("C1147").Value = '1730'
$ws.Range("D1147").Value = 'CLH'
$ws.Range("E1147").Value = 'I'
$ws.Range("F1147").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1148
$ws.Range("A1148").Value = 'AV Shutdown'
$ws.Range("B1148").Value = 42755
$ws.Range("C1148").Value = '1730'
$ws.Range("D1148").Value = 'LAS'
$ws.Range("E1148").Value = 'A'
$ws.Range("F1148").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1149
$ws.Range("A1149").Value = 'AV Shutdown'
$ws.Range("B1149").Value = 42755
$ws.Range("C1149").Value = '1630'
$ws.Range("D1149").Value = 'LAS'
$ws.Range("E1149").Value = 'C'
$ws.Range("F1149").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1150
$ws.Range("A1150").Value = 'AV Shutdown'
$ws.Range("B1150").Value = 42755
$ws.Range("C1150").Value = '1600'
$ws.Range("D1150").Value = 'LSB'
$ws.Range("E1150").Value = '106'
$ws.Range("F1150").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# Row 1151
$ws.Range("A1151").Value = 'AV Shutdown'
$ws.Range("B1151").Value = 42755
$ws.Range("C1151").Value = '1630'
$ws.Range("D1151").Value = 'LSB'
$ws.Range("E1151").Value = '105'
$ws.Range("F1151").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

$ws.Application.CutCopyMode = $false

# Update sheet view to match final cursor position/selection
$ws.Range("F1155").Select()